# Update Leve profit sheets with refreshed market-board pricing data
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 316.1111
$ws.Range("I28").Value = 316.1111
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 316.1111
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 168.8889
$ws.Range("N28").Value = ""
$ws.Range("H32").Value = 9011015
$ws.Range("I32").Value = 33334908
$ws.Range("J32").Value = 2165.1482
$ws.Range("K32").Value = 33334908
$ws.Range("L32").Value = 2165.1482
$ws.Range("M32").Value = -33334582
$ws.Range("N32").Value = -2817.1482
$ws.Range("H106").Value = 67548.5
$ws.Range("I106").Value = 2902.5
$ws.Range("J106").Value = 80477.7
$ws.Range("K106").Value = 2902.5
$ws.Range("L106").Value = 80477.7
$ws.Range("M106").Value = -2271.5
$ws.Range("N106").Value = -81739.7
$ws.Range("H132").Value = 25140.512
$ws.Range("I132").Value = 3455.9092
$ws.Range("J132").Value = 114589.5
$ws.Range("K132").Value = 10367.7276
$ws.Range("L132").Value = 343768.5
$ws.Range("M132").Value = -7837.7276
$ws.Range("N132").Value = -348828.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3555.5557
$ws.Range("I2").Value = 2750
$ws.Range("J2").Value = 10000
$ws.Range("K2").Value = 2750
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = -2637
$ws.Range("N2").Value = -10226
$ws.Range("H4").Value = 466.33334
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 466.33334
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 466.33334
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = -698.33334
$ws.Range("H5").Value = 335.9
$ws.Range("I5").Value = 219.625
$ws.Range("J5").Value = 801
$ws.Range("K5").Value = 219.625
$ws.Range("L5").Value = 801
$ws.Range("M5").Value = -107.625
$ws.Range("N5").Value = -1025
$ws.Range("H26").Value = 7266.3335
$ws.Range("I26").Value = 7266.3335
$ws.Range("K26").Value = 7266.3335
$ws.Range("M26").Value = -6936.3335
$ws.Range("H39").Value = 70018
$ws.Range("J39").Value = 70018
$ws.Range("L39").Value = 70018
$ws.Range("N39").Value = -71058
$ws.Range("H74").Value = 1653.7
$ws.Range("J74").Value = 2696.625
$ws.Range("L74").Value = 2696.625
$ws.Range("N74").Value = -4444.625
$ws.Range("H77").Value = 1653.7
$ws.Range("J77").Value = 2696.625
$ws.Range("L77").Value = 13483.125
$ws.Range("N77").Value = -22219.125
$ws.Range("H116").Value = 3555.5557
$ws.Range("I116").Value = 2750
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 2750
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = -456
$ws.Range("N116").Value = -14588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3555.5557
$ws.Range("I3").Value = 2750
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 2750
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = -2636
$ws.Range("N3").Value = -10228
$ws.Range("H4").Value = 335.9
$ws.Range("I4").Value = 219.625
$ws.Range("J4").Value = 801
$ws.Range("K4").Value = 219.625
$ws.Range("L4").Value = 801
$ws.Range("M4").Value = -104.625
$ws.Range("N4").Value = -1031
$ws.Range("H22").Value = 220
$ws.Range("I22").Value = 220
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 220
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -47
$ws.Range("N22").Value = ""
$ws.Range("H105").Value = 4380.9165
$ws.Range("I105").Value = 2497.5
$ws.Range("J105").Value = 5322.625
$ws.Range("K105").Value = 2497.5
$ws.Range("L105").Value = 5322.625
$ws.Range("M105").Value = -750.5
$ws.Range("N105").Value = -8816.625
$ws.Range("H137").Value = 64221.715
$ws.Range("J137").Value = 64221.715
$ws.Range("L137").Value = 64221.715
$ws.Range("N137").Value = -74421.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1868.6666
$ws.Range("I16").Value = 1433
$ws.Range("K16").Value = 1433
$ws.Range("M16").Value = -1146
$ws.Range("H31").Value = 4978500
$ws.Range("I31").Value = 1486.3871
$ws.Range("J31").Value = 9264261
$ws.Range("K31").Value = 1486.3871
$ws.Range("L31").Value = 9264261
$ws.Range("M31").Value = -1191.3871
$ws.Range("N31").Value = -9264851
$ws.Range("H34").Value = 4978500
$ws.Range("I34").Value = 1486.3871
$ws.Range("J34").Value = 9264261
$ws.Range("K34").Value = 1486.3871
$ws.Range("L34").Value = 9264261
$ws.Range("M34").Value = -1284.3871
$ws.Range("N34").Value = -9264665
$ws.Range("H35").Value = 935.7143
$ws.Range("I35").Value = 745.36365
$ws.Range("J35").Value = 1633.6666
$ws.Range("K35").Value = 745.36365
$ws.Range("L35").Value = 1633.6666
$ws.Range("M35").Value = -451.36365
$ws.Range("N35").Value = -2221.6666
$ws.Range("H99").Value = 3168.6667
$ws.Range("I99").Value = 3168.6667
$ws.Range("K99").Value = 3168.6667
$ws.Range("M99").Value = -1670.6667
$ws.Range("H113").Value = 1868.6666
$ws.Range("I113").Value = 1433
$ws.Range("K113").Value = 1433
$ws.Range("M113").Value = 737
$ws.Range("H126").Value = 3168.6667
$ws.Range("I126").Value = 3168.6667
$ws.Range("K126").Value = 9506.000100000001
$ws.Range("M126").Value = -7036.000100000001
$ws.Range("H134").Value = 2003069.4
$ws.Range("I134").Value = 2496.3333
$ws.Range("J134").Value = 3503499.2
$ws.Range("K134").Value = 7488.999899999999
$ws.Range("L134").Value = 10510497.6
$ws.Range("M134").Value = -4953.999899999999
$ws.Range("N134").Value = -10515567.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1281.7073
$ws.Range("I68").Value = 1116.9
$ws.Range("J68").Value = 1334.871
$ws.Range("K68").Value = 3350.7
$ws.Range("L68").Value = 4004.613
$ws.Range("M68").Value = -2539.7
$ws.Range("N68").Value = -5626.613
$ws.Range("H71").Value = 1281.7073
$ws.Range("I71").Value = 1116.9
$ws.Range("J71").Value = 1334.871
$ws.Range("K71").Value = 10052.1
$ws.Range("L71").Value = 12013.839
$ws.Range("M71").Value = -5996.1
$ws.Range("N71").Value = -20125.839
$ws.Range("H107").Value = 2719.2341
$ws.Range("I107").Value = 3722.8965
$ws.Range("J107").Value = 2271.446
$ws.Range("K107").Value = 11168.6895
$ws.Range("L107").Value = 6814.338
$ws.Range("M107").Value = -9248.6895
$ws.Range("N107").Value = -10654.338
$ws.Range("H140").Value = 3275.4167
$ws.Range("I140").Value = 857.25
$ws.Range("J140").Value = 4484.5
$ws.Range("K140").Value = 2571.75
$ws.Range("L140").Value = 13453.5
$ws.Range("M140").Value = 2608.25
$ws.Range("N140").Value = -23813.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 15443.272
$ws.Range("I126").Value = 39039
$ws.Range("J126").Value = 1960
$ws.Range("K126").Value = 117117
$ws.Range("L126").Value = 5880
$ws.Range("M126").Value = -114647
$ws.Range("N126").Value = -10820
$ws.Range("H140").Value = 37251.9
$ws.Range("J140").Value = 37251.9
$ws.Range("L140").Value = 37251.9
$ws.Range("N140").Value = -47611.9
$ws.Range("H141").Value = 70105.75
$ws.Range("J141").Value = 70105.75
$ws.Range("L141").Value = 70105.75
$ws.Range("N141").Value = -80465.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 39994.332
$ws.Range("J140").Value = 39994.332
$ws.Range("L140").Value = 39994.332
$ws.Range("N140").Value = -50354.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 38294.848
$ws.Range("I29").Value = 2666.6667
$ws.Range("J29").Value = 48983.3
$ws.Range("K29").Value = 2666.6667
$ws.Range("L29").Value = 48983.3
$ws.Range("M29").Value = -2376.6667
$ws.Range("N29").Value = -49563.3
$ws.Range("H140").Value = 34188.152
$ws.Range("J140").Value = 34188.152
$ws.Range("L140").Value = 34188.152
$ws.Range("N140").Value = -44548.152
$ws.Range("H141").Value = 39765.4
$ws.Range("J141").Value = 39765.4
$ws.Range("L141").Value = 39765.4
$ws.Range("N141").Value = -50125.4
